# Slide 9 ("Recommendations") content placeholder: merge the split runs in
# paragraphs 2 and 3 back into single runs (same visible text, but the
# run-level XML had been split into two <a:r> elements that should really
# be one).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# --- Paragraph 2: "Choose aircraft with a record of mechanical " + "reliability."
# Merging keeps the FIRST run's formatting (it already carries dirty="0"),
# so simply overwrite the whole paragraph's text with the concatenation.
# (Use Characters(start,len) rather than the Paragraphs() range itself --
# re-assigning .Text through Paragraphs() leaves the original run split in
# this runtime, whereas a plain Characters() sub-range correctly collapses
# it down to a single run.)
$para2 = $tr.Paragraphs(2, 1)
$para2Start = $para2.Start
$para2Length = $para2.Length
$para2Range = $tr.Characters($para2Start, $para2Length)
$para2Range.Text = "Choose aircraft with a record of mechanical reliability."

# --- Paragraph 3: "Implement " + "strict weather-related flight operation protocols."
# Merging keeps the SECOND run's formatting (the one with dirty="0"), so
# remove the first run's text, then let the remaining run absorb the prefix.
# (Capture Start/Length as plain numbers first -- the Paragraphs()/Characters()
# ranges are "live" and their Length shrinks as soon as text is deleted.)
$para3 = $tr.Paragraphs(3, 1)
$para3Start = $para3.Start
$para3Length = $para3.Length
$firstRunLen = "Implement ".Length

$firstRun = $tr.Characters($para3Start, $firstRunLen)
$firstRun.Text = ""

$remaining = $tr.Characters($para3Start, $para3Length - $firstRunLen)
$remaining.Text = "Implement strict weather-related flight operation protocols."
